$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Update D3: Status changes from "In Progress" to "Completed"
$ws.Range("D3").Value = "Completed"

# Update E3: add new comment string
$ws.Range("E3").Value = "05-04-Supreet - Formatted and uploaded into GIT - 05-09"

# Update selection to D4
$ws.Activate()
$ws.Range("D4").Select()
